# This workbook's rows 6-14 each hold one species observation record.
# The edit re-assigns which observation record sits on which row: the
# content of several columns (A, B, D, E, F, G, H, Q, R, AC) moves between
# rows 6-14 following a permutation, while the rest of each row's columns
# (which are identical across all these rows anyway) stay put.
#
# Row (new) <- Row (old, source of the record that lands there)
#   6  <- 11
#   7  <- 10
#   8  <- 7
#   9  <- 13
#  10  <- 14
#  11  <- 9
#  12  <- 6
#  13  <- 12
#  14  <- 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (1-based) that travel with each observation record.
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18, 29)   # A B D E F G H Q R AC

# 1) Snapshot the current (pre-edit) values for every tracked column,
#    for every affected row, BEFORE any writes happen.
$snapshot = @{}
for ($r = 6; $r -le 14; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Mapping: for each destination row, which row's snapshot to copy from.
$mapping = @{
    6  = 11
    7  = 10
    8  = 7
    9  = 13
    10 = 14
    11 = 9
    12 = 6
    13 = 12
    14 = 8
}

# 3) Write the snapshot values back out according to the mapping.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
